$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$json = '{"quiz_questions":[{"id":"1","question":"Sparing in words, using very few words","possible_answers":"Loquacious,Laconic,Integral,Judicious","correct_answer":"2"},{"id":"2","question":"Freedom from punishment","possible_answers":"Impunity,Indolent,Jaded,Hedonist","correct_answer":"1"},{"id":"3","question":"An instrument for measuring","possible_answers":"Gallant ,Gauge ,Gamester ,Gastric ","correct_answer":"2"},{"id":"4","question":"Inflammation of the stomach  ","possible_answers":"Insurgent ,Impudence ,Garrulous ,Gastritis ","correct_answer":"4"},{"id”:”5”,”question":"unoriginal","possible_answers":"Imperious ,Impudence ,Hackneyed ,Inane ","correct_answer":"3"},{"id”:”6”,”question":"To draw principle inferences","possible_answers":"Gaseous ,Generalize ,Garrote ,Garrison ","correct_answer":"2"},{"id”:”7”,”question":"To imbue with life or animation ","possible_answers":"Galvanize ,Garnish ,Garrulous ,Gaseous ","correct_answer":"1"},{"id”:”8”,”question":"Decorate or embellish ","possible_answers":"Gamut ,Genealogist ,Gambol ,Garnish ","correct_answer":"4"},{"id”:”9”,”question":"Given to constant trivial talking","possible_answers":"Gaseous ,Gendarme ,Garrulous ,Genitive ","correct_answer":"3"},{"id”:”10”,”question":"To risk money or other possession on an event","possible_answers":"Gallant ,Gaily ,Gamble ,Gastronomy ","correct_answer":"3"}]}'

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = $json

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = $json

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = $json
